$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Message literals (order below mirrors the shared-string insertion order of the target workbook)
$succes = 'Succès'
$demoDesc = 'This is a demo for NORAUI (Non-Regression Automation for User Interfaces).'
$apiJson = '{"total_count":218,"incomplete_results":false,"items":[{"login":"ndeloof","id":132757,"node_id":"MDQ6VXNlcjEzMjc1Nw==","avatar_url":"https://avatars0.githubusercontent.com/u/132757?v=4","gravatar_id":"","url":"https://api.github.com/users/ndeloof","html_url":"https://github.com/ndeloof","followers_url":"https://api.github.com/users/ndeloof/followers","following_url":"https://api.github.com/users/ndeloof/following{/other_user}","gists_url":"https://api.github.com/users/ndeloof/gists{/gist_id}","starred_url":"https://api.github.com/users/ndeloof/starred{/owner}{/repo}","subscriptions_url":"https://api.github.com/users/ndeloof/subscriptions","organizations_url":"https://api.github.com/users/ndeloof/orgs","repos_url":"https://api.github.com/users/ndeloof/repos","events_url":"https://api.github.com/users/ndeloof/events{/privacy}","received_events_url":"https://api.github.com/users/ndeloof/received_events","type":"User","site_admin":false,"score":1.0},{"login":"aurelhubert","id":4895410,"node_id":"MDQ6VXNlcjQ4OTU0MTA=","avatar_url":"https://avatars0.githubusercontent.com/u/4895410?v=4","gravatar_id":"","url":"https://api.github.com/users/aurelhubert","html_url":"https://github.com/aurelhubert","followers_url":"https://api.github.com/users/aurelhubert/followers","following_url":"https://api.github.com/users/aurelhubert/following{/other_user}","gists_url":"https://api.github.com/users/aurelhubert/gists{/gist_id}","starred_url":"https://api.github.com/users/aurelhubert/starred{/owner}{/repo}","subscriptions_url":"https://api.github.com/users/aurelhubert/subscriptions","organizations_url":"https://api.github.com/users/aurelhubert/orgs","repos_url":"https://api.github.com/users/aurelhubert/repos","events_url":"https://api.github.com/users/aurelhubert/events{/privacy}","received_events_url":"https://api.github.com/users/aurelhubert/received_events","type":"User","site_admin":false,"score":1.0},{"login":"pcmanus","id":479957,"node_id":"MDQ6VXNlcjQ3OTk1Nw==","avatar_url":"https://avatars3.githubusercontent.com/u/479957?v=4","gravatar_id":"","url":"https://api.github.com/users/pcmanus","html_url":"https://github.com/pcmanus","followers_url":"https://api.github.com/users/pcmanus/followers","following_url":"https://api.github.com/users/pcmanus/following{/other_user}","gists_url":"https://api.github.com/users/pcmanus/gists{/gist_id}","starred_url":"https://api.github.com/users/pcmanus/starred{/owner}{/repo}","subscriptions_url":"https://api.github.com/users/pcmanus/subscriptions","organizations_url":"https://api.github.com/users/pcmanus/orgs","repos_url":"https://api.github.com/users/pcmanus/repos","events_url":"https://api.github.com/users/pcmanus/events{/privacy}","received_events_url":"https://api.github.com/users/pcmanus/received_events","type":"User","site_admin":false,"score":1.0},{"login":"deruelle","id":54279,"node_id":"MDQ6VXNlcjU0Mjc5","avatar_url":"https://avatars3.githubusercontent.com/u/54279?v=4","gravatar_id":"","url":"https://api.github.com/users/deruelle","html_url":"https://github.com/deruelle","followers_url":"https://api.github.com/users/deruelle/followers","following_url":"https://api.github.com/users/deruelle/following{/other_user}","gists_url":"https://api.github.com/users/deruelle/gists{/gist_id}","starred_url":"https://api.github.com/users/deruelle/starred{/owner}{/repo}","subscriptions_url":"https://api.github.com/users/deruelle/subscriptions","organizations_url":"https://api.github.com/users/deruelle/orgs","repos_url":"https://api.github.com/users/deruelle/repos","events_url":"https://api.github.com/users/deruelle/events{/privacy}","received_events_url":"https://api.github.com/users/deruelle/received_events","type":"User","site_admin":false,"score":1.0},{"login":"barais","id":598311,"node_id":"MDQ6VXNlcjU5ODMxMQ==","avatar_url":"https://avatars1.githubusercontent.com/u/598311?v=4","gravatar_id":"","url":"https://api.github.com/users/barais","html_url":"https://github.com/barais","followers_url":"https://api.github.com/users/barais/followers","following_url":"https://api.github.com/users/barais/following{/other_user}","gists_url":"https://api.github.com/users/barais/gists{/gist_id}","starred_url":"https://api.github.com/users/barais/starred{/owner}{/repo}","subscriptions_url":"https://api.github.com/users/barais/subscriptions","organizations_url":"https://api.github.com/users/barais/orgs","repos_url":"https://api.github.com/users/barais/repos","events_url":"https://api.github.com/users/barais/events{/privacy}","received_events_url":"https://api.github.com/users/barais/received_events","type":"User","site_admin":false,"score":1.0},{"login":"ybonnel","id":428232,"node_id":"MDQ6VXNlcjQyODIzMg==","avatar_url":"https://avatars1.githubusercontent.com/u/428232?v=4","gravatar_id":"","url":"https://api.github.com/users/ybonnel","html_url":"https://github.com/ybonnel","followers_url":"https://api.github.com/users/ybonnel/followers","following_url":"https://api.github.com/users/ybonnel/following{/other_user}","gists_url":"https://api.github.com/users/ybonnel/gists{/gist_id}","starred_url":"https://api.github.com/users/ybonnel/starred{/owner}{/repo}","subscriptions_url":"https://api.github.com/users/ybonnel/subscriptions","organizations_url":"https://api.github.com/users/ybonnel/orgs","repos_url":"https://api.github.com/users/ybonnel/repos","events_url":"https://api.github.com/users/ybonnel/events{/privacy}","received_events_url":"https://api.github.com/users/ybonnel/received_events","type":"User","site_admin":false,"score":1.0},{"login":"ganfra","id":8754928,"node_id":"MDQ6VXNlcjg3NTQ5Mjg=","avatar_url":"https://avatars3.githubusercontent.com/u/8754928?v=4","gravatar_id":"","url":"https://api.github.com/users/ganfra","html_url":"https://github.com/ganfra","followers_url":"https://api.github.com/users/ganfra/followers","following_url":"https://api.github.com/users/ganfra/following{/other_user}","gists_url":"https://api.github.com/users/ganfra/gists{/gist_id}","starred_url":"https://api.github.com/users/ganfra/starred{/owner}{/repo}","subscriptions_url":"https://api.github.com/users/ganfra/subscriptions","organizations_url":"https://api.github.com/users/ganfra/orgs","repos_url":"https://api.github.com/users/ganfra/repos","events_url":"https://api.github.com/users/ganfra/events{/privacy}","received_events_url":"https://api.github.com/users/ganfra/received_events","type":"User","site_admin":false,"score":1.0},{"login":"Jean-Emile","id":1098811,"node_id":"MDQ6VXNlcjEwOTg4MTE=","avatar_url":"https://avatars3.githubusercontent.com/u/1098811?v=4","gravatar_id":"","url":"https://api.github.com/users/Jean-Emile","html_url":"https://github.com/Jean-Emile","followers_url":"https://api.github.com/users/Jean-Emile/followers","following_url":"https://api.github.com/users/Jean-Emile/following{/other_user}","gists_url":"https://api.github.com/users/Jean-Emile/gists{/gist_id}","starred_url":"https://api.github.com/users/Jean-Emile/starred{/owner}{/repo}","subscriptions_url":"https://api.github.com/users/Jean-Emile/subscriptions","organizations_url":"https://api.github.com/users/Jean-Emile/orgs","repos_url":"https://api.github.com/users/Jean-Emile/repos","events_url":"https://api.github.com/users/Jean-Emile/events{/privacy}","received_events_url":"https://api.github.com/users/Jean-Emile/received_events","type":"User","site_admin":false,"score":1.0},{"login":"arnobl","id":1052622,"node_id":"MDQ6VXNlcjEwNTI2MjI=","avatar_url":"https://avatars0.githubusercontent.com/u/1052622?v=4","gravatar_id":"","url":"https://api.github.com/users/arnobl","html_url":"https://github.com/arnobl","followers_url":"https://api.github.com/users/arnobl/followers","following_url":"https://api.github.com/users/arnobl/following{/other_user}","gists_url":"https://api.github.com/users/arnobl/gists{/gist_id}","starred_url":"https://api.github.com/users/arnobl/starred{/owner}{/repo}","subscriptions_url":"https://api.github.com/users/arnobl/subscriptions","organizations_url":"https://api.github.com/users/arnobl/orgs","repos_url":"https://api.github.com/users/arnobl/repos","events_url":"https://api.github.com/users/arnobl/events{/privacy}","received_events_url":"https://api.github.com/users/arnobl/received_events","type":"User","site_admin":false,"score":1.0},{"login":"acherm","id":530573,"node_id":"MDQ6VXNlcjUzMDU3Mw==","avatar_url":"https://avatars1.githubusercontent.com/u/530573?v=4","gravatar_id":"","url":"https://api.github.com/users/acherm","html_url":"https://github.com/acherm","followers_url":"https://api.github.com/users/acherm/followers","following_url":"https://api.github.com/users/acherm/following{/other_user}","gists_url":"https://api.github.com/users/acherm/gists{/gist_id}","starred_url":"https://api.github.com/users/acherm/starred{/owner}{/repo}","subscriptions_url":"https://api.github.com/users/acherm/subscriptions","organizations_url":"https://api.github.com/users/acherm/orgs","repos_url":"https://api.github.com/users/acherm/repos","events_url":"https://api.github.com/users/acherm/events{/privacy}","received_events_url":"https://api.github.com/users/acherm/received_events","type":"User","site_admin":false,"score":1.0}]}'
$failParis = 'Échec : The city is Paris!!'
$failInputSelect = 'Échec : /!\ Saisie « Input Select field » dans demo. /!\'
$failNoExist = 'Échec : /!\ Accès à l''action « no exist element » dans demo. /!\'
$failVilleVide = 'Échec : /!\ La donnée « la ville » fournie ne peut pas être vide. /!\'
$failElementsVide = 'Échec : /!\ La donnée « élément » fournie ne peut pas être vide. /!\, /!\ La donnée « élément2 » fournie ne peut pas être vide. /!\.'
$failInputTextField = 'Échec : /!\ Absence « -input_text_field » dans demo. /!\'

$greenColorIndex = 10  # renders as a dark green font (success)
$redColorIndex = 3     # renders as a red font (failure)

# --- Row 2 : Rennes scenario -> success ---
$ws.Range("H2").Value = $succes
$ws.Range("H2").Font.ColorIndex = $greenColorIndex

# --- Row 9 : demo description (REST payload column reused to describe the demo) ---
$g9 = $ws.Range("G9")
$g9.ClearFormats()
$g9.Value = $demoDesc

# --- Row 2 : REST API payload returned by the call ---
$g2 = $ws.Range("G2")
$g2.ClearFormats()
$g2.Value = $apiJson

# --- Row 3 : Paris scenario -> failure (unexpected city) ---
$ws.Range("H3").Value = $failParis
$ws.Range("H3").Font.ColorIndex = $redColorIndex

# --- Row 4 : failure entering the input select field ---
$ws.Range("H4").Value = $failInputSelect
$ws.Range("H4").Font.ColorIndex = $redColorIndex

# --- Row 6 : failure accessing an unknown element ---
$ws.Range("H6").Value = $failNoExist
$ws.Range("H6").Font.ColorIndex = $redColorIndex

# --- Row 7 : failure, empty ville data ---
$ws.Range("H7").Value = $failVilleVide
$ws.Range("H7").Font.ColorIndex = $redColorIndex

# --- Row 8 : failure, empty élément/élément2 data ---
$ws.Range("H8").Value = $failElementsVide
$ws.Range("H8").Font.ColorIndex = $redColorIndex

# --- Row 9 : failure, missing input text field ---
$ws.Range("H9").Value = $failInputTextField
$ws.Range("H9").Font.ColorIndex = $redColorIndex

# --- Row 5 : Rennes scenario -> success (reuses strings already interned above) ---
$ws.Range("H5").Value = $succes
$ws.Range("H5").Font.ColorIndex = $greenColorIndex
$g5 = $ws.Range("G5")
$g5.ClearFormats()
$g5.Value = $apiJson

Write-Host "edit applied"
